$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new permit site ("Karasjok") was added to the source CSV between the
# existing "Caskin-jeaggi / Caskinjohka" (row 5) and "Fahttevarleaksa" (old
# row 6) records. Re-running the query effectively inserts one new row and
# shifts everything below it down by one.
$ws.Rows("6:6").Insert()

# Populate the newly inserted row with the Karasjok site data.
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Karasjok"
$ws.Range("C6").Value = 25.386399999999998
$ws.Range("D6").Value = 69.444800000000001
$ws.Range("E6").Value = 436777.59009725502
$ws.Range("F6").Value = 7704794.2822894901
$ws.Range("G6").Value = 35
$ws.Range("H6").Value = "Karasjok"

# The sequential "Site ID" numbers (column A) for every row below the
# inserted one simply advance by one.
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10

# Column H ("Kommune") re-fits its width now that the data set has changed.
$ws.Columns("H:H").ColumnWidth = 8.764322916666666

# The query table now spans one additional row, so the workbook-level
# defined name that mirrors its extent needs to grow accordingly.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*permits_sept2021_all_sites") {
        $n.RefersTo = "=Sheet1!`$B`$1:`$H`$11"
    }
}

# Reflect the refreshed query table's new extent in the sheet's selection.
$ws.Range("A1:H11").Select()
